$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53; this shifts the existing rows 53-115
# down to 54-116 (carrying their values/styles along), matching the
# target diff where every row from 53 downward is the prior row's data.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new weekly record.
$ws.Range("A53").Value = 5
$ws.Range("B53").Value = "Macroferia Regional de Talca"
$ws.Range("C53").Value = "Maule"
$ws.Range("D53").Value = 45159
$ws.Range("E53").Value = 7
$ws.Range("F53").Value = 100112040
$ws.Range("G53").Value = "Cilantro"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 300
$ws.Range("K53").Value = 800
$ws.Range("L53").Value = 800
$ws.Range("M53").Value = 800
$ws.Range("N53").Value = "`$/caja 36 atados"
$ws.Range("O53").Value = "Región Metropolitana"
$ws.Range("P53").Value = 22
$ws.Range("Q53").Value = 36
$ws.Range("R53").Value = "Hortaliza"
